$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Build/knit date (appears on title page and in the session-info "date" line;
# ReplaceAll below updates both occurrences in one pass)
Replace-Text "2019-10-12" "2019-10-28"

# Pathway table 1 (Metabolic pathways) - Beta-Alanine Metabolism row
Replace-Text "0.005" "0.006"
Replace-Text "0.279" "0.226"

# Histidine Metabolism row
Replace-Text "0.008" "0.006"
Replace-Text "0.279" "0.226"

# Glycine and Serine Metabolism row
Replace-Text "0.042" "0.047"
Replace-Text "0.849" "0.836"

# Pathway table 2 - Sphingolipid Metabolism row
Replace-Text "0.051" "0.068"
Replace-Text "0.988" "0.989"

# Phospholipid Biosynthesis row
Replace-Text "0.090" "0.105"
Replace-Text "0.988" "0.989"

# Session info package versions
Replace-Text "ezlimma      * 0.2.3.9028 2019-10-11 [1] local         " "ezlimma      * 0.2.3.9029 2019-10-22 [1] local         "
Replace-Text "ezlimmaplot  * 0.0.1.9016 2019-10-07 [1] local         " "ezlimmaplot  * 0.0.1.9017 2019-10-24 [1] local         "
Replace-Text "Hitman       * 0.0.0.9006 2019-10-06 [1] local         " "Hitman       * 0.0.0.9006 2019-10-21 [1] local         "
